$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same header formatting (bold, centered, bordered) used by the
# existing header row to the new columns, by copying the style from AC1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Add header labels in row 1 for the new team-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill data rows 2-56 with the team's record (same W/L/T for every player row)
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 76   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 86   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
